$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-11 (they are no longer present in the new scenario)
$ws.Range("A8:H11").Delete()

# Row 2: new opening system utterance; the old user-example and condition are removed
$ws.Range("C2").Value = 'こんにちは。私はチャットボットです。気軽にお話しましょう。料理は良くする方ですか？ '
$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3 becomes "state1" with a new generation prompt & a boredom-check condition,
# looping back to state1
$nl = "`n"
$ws.Range("B3").Value = 'state1'
$ws.Range("C3").Value = '$あなたの次の発話を30文字以内で生成してください。$' + $nl
$ws.Range("F3").Value = '$ユーザが飽きているかどうか判断してください$'
$ws.Range("H3").Value = 'state1'

# Row 4 stays state1 -> state2, now gated on a turn-count condition, no system utterance
$ws.Range("C4").ClearContents()
$ws.Range("F4").Value = 'TS>5'
$ws.Range("H4").Value = 'state2'

# Row 5 becomes #error (previously state2)
$ws.Range("B5").Value = '#error'
$ws.Range("C5").Value = '申し訳ありません。内部エラーがおきてしまいました。今日はありがとうございました。'
$ws.Range("D5").ClearContents()
$ws.Range("H5").ClearContents()

# Row 6 becomes #final_state1 (previously state3)
$ws.Range("B6").Value = '#final_state1'
$ws.Range("C6").Value = '$それまでの発話に続けて、対話を終わらせる発話を30文字以内で生成してください。$'
$ws.Range("F6").ClearContents()
$ws.Range("H6").ClearContents()

# Row 7 becomes state2 with the new persona-based chit-chat prompt
$promptLines = @(
    '$$$',
    '# 状況',
    '{situation}',
    '# あなたのペルソナ',
    '{persona}',
    '# 今までの対話',
    '{dialogue_history}',
    '# タスク',
    '料理以外の話に関して雑談をするために、あなたの次の発話を30文字以内で生成してください。',
    '$$$ '
)
$c7 = ($promptLines -join $nl) + $nl
$ws.Range("B7").Value = 'state2'
$ws.Range("C7").Value = $c7
$ws.Range("F7").Value = 'TT>10'
$ws.Range("H7").Value = '#final_state1'

Write-Host "Done. Used range: $($ws.UsedRange.Address())"
